$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1824
$ws.Range("J3").Value = 1893
$ws.Range("I4").Value = 1758
$ws.Range("J4").Value = 427
$ws.Range("J5").Value = 131
$ws.Range("J6").Value = 2429
$ws.Range("I7").Value = 26205
$ws.Range("J7").Value = 6704

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J7").Value = 188
$ws.Range("J8").Value = 409
$ws.Range("J15").Value = 86
$ws.Range("J18").Value = 79
$ws.Range("J19").Value = 228
$ws.Range("J20").Value = 141
$ws.Range("J22").Value = 13
$ws.Range("J29").Value = 376
$ws.Range("J33").Value = 277
$ws.Range("J37").Value = 223
$ws.Range("J42").Value = 266
$ws.Range("J43").Value = 69
$ws.Range("J48").Value = 56
$ws.Range("J49").Value = 39
$ws.Range("J51").Value = 93
$ws.Range("J52").Value = 155
$ws.Range("J54").Value = 135
$ws.Range("J57").Value = 35
$ws.Range("J58").Value = 3
$ws.Range("J62").Value = 2
$ws.Range("J63").Value = 32
$ws.Range("J65").Value = 178
$ws.Range("J67").Value = 245
$ws.Range("J70").Value = 12
$ws.Range("J77").Value = 48
$ws.Range("J83").Value = 160
$ws.Range("J85").Value = 311
$ws.Range("J87").Value = 26
$ws.Range("J90").Value = 79
$ws.Range("J92").Value = 24
$ws.Range("J93").Value = 30
$ws.Range("J95").Value = 100
$ws.Range("I96").Value = 305
$ws.Range("J99").Value = 87
$ws.Range("I101").Value = 26205
$ws.Range("J101").Value = 6704

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 127
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 134
$ws.Range("J3").Value = 136
$ws.Range("J7").Value = 409

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 59
$ws.Range("J3").Value = 61
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 305

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 77
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 98
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 245

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 53
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 55
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 70
$ws.Range("J3").Value = 79
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 277

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 67
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J3").Value = 134
$ws.Range("J5").Value = 14
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 376

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 63
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 8
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 55
$ws.Range("J6").Value = 138
$ws.Range("J7").Value = 266

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 141

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J3").Value = 2
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J4").Value = 16
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J2").Value = 7
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 3

$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Range("J2").Value = 1
$ws.Range("J6").Value = 2
